# Updated cryptos list on Sat Jul 29 07:23:33 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns; a couple of rank swaps too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text even when the literal looks numeric
    # (e.g. "0.7125", "241.70", "29.350.06"), then drop back to the
    # default "Normal" style so no stray number-format sticks around.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "29.350.06"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.873.76"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - XRP
Set-TextValue $ws.Range("D5") "0.7125"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6 - BNB
Set-TextValue $ws.Range("D6") "241.70"
$ws.Range("E6").Value = "  +0.23%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3115"
$ws.Range("E8").Value = "  +1.03%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.07775"
$ws.Range("E9").Value = "  +1.47%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "25.12"
$ws.Range("E10").Value = "  +1.51%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.11%  "

# Row 12 - WrappedEther
Set-TextValue $ws.Range("D12") "1.874.45"
$ws.Range("E12").Value = "  +0.62%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("D13") "5.233"
$ws.Range("E13").Value = "  +0.83%  "

# Row 14 - Polygon
Set-TextValue $ws.Range("D14") "0.7125"
$ws.Range("E14").Value = "  +0.29%  "

# Row 15 - Litecoin
Set-TextValue $ws.Range("D15") "91.12"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "29.357.17"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17 - Uniswap
Set-TextValue $ws.Range("D17") "6.094"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.000008244"
$ws.Range("E18").Value = "  +5.46%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "240.36"
$ws.Range("E19").Value = "  -1.04%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  +0.94%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D21") "2.122.04"
$ws.Range("E21").Value = "  +0.29%  "

# Row 22 - Dai
Set-TextValue $ws.Range("D22") "0.9998"

# Row 23 - Chainlink
Set-TextValue $ws.Range("D23") "7.767"
$ws.Range("E23").Value = "  -1.47%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  +0.00%  "

# Row 25 - Stellar
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "163.08"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.054"
$ws.Range("E27").Value = "  +1.48%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.40%  "

# Row 29 - PancakeSwap
Set-TextValue $ws.Range("D29") "1.511"
$ws.Range("E29").Value = "  +0.80%  "

# Row 30 - Filecoin
$ws.Range("E30").Value = "  +0.38%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D31") "4.324"
$ws.Range("E31").Value = "  +2.54%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  -1.92%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.05286"
$ws.Range("E33").Value = "  +3.16%  "

# Row 34 - LidoDAOToken
Set-TextValue $ws.Range("D34") "1.937"
$ws.Range("E34").Value = "  +1.43%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +1.32%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("D36") "0.7404"
$ws.Range("E36").Value = "  -8.53%  "

# Row 37 - HuobiToken
Set-TextValue $ws.Range("D37") "2.701"
$ws.Range("E37").Value = "  +0.77%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +1.84%  "

# Row 39 - Maker
Set-TextValue $ws.Range("D39") "1.216.87"
$ws.Range("E39").Value = "  +4.19%  "

# Row 40 - MXToken
Set-TextValue $ws.Range("D40") "2.730"
$ws.Range("E40").Value = "  +1.20%  "

# Row 41 - FraxShare
Set-TextValue $ws.Range("D41") "6.544"
$ws.Range("E41").Value = "  +5.93%  "

# Row 42 - Quant
Set-TextValue $ws.Range("D42") "110.96"
$ws.Range("E42").Value = "  +8.69%  "

# Row 43 - now Aave (was TrustWalletToken)
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D43") "72.92"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44 - now TrustWalletToken (was Aave)
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D44") "0.8874"
$ws.Range("E44").Value = "  -0.42%  "

# Row 46 - RocketPoolETH
Set-TextValue $ws.Range("D46") "2.021.21"
$ws.Range("E46").Value = "  +0.40%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +1.96%  "

# Row 48 - Mantle
Set-TextValue $ws.Range("D48") "0.5211"
$ws.Range("E48").Value = "  +0.66%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +2.42%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "9.404"
$ws.Range("E50").Value = "  +1.44%  "

# Row 51 - TheSandbox
Set-TextValue $ws.Range("D51") "0.4318"
$ws.Range("E51").Value = "  +1.10%  "
